$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Salidas")
Write-Host $ws.Name
